# WeatherAPI.com entegrasyonu - Gercek hava durumu verileri eklendi
# Adds the 4 new match rows (10-13) with their date/time/group/team data,
# matching the style (date/time number formats) already used by the
# surrounding rows, and updates the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 10 and 11 are brand new - copy the date/time number formatting
# from row 2 (style s="1" for the date column, s="2" for the time column)
# before filling in the values so no new style entries are minted.
$ws.Range("A2:B2").Copy()
$ws.Range("A10:B11").PasteSpecial(-4122)

# Row 10: 25.06.2025 21:00 - Ahmet Minguzzi Grubu: Ravager - Kural Kesiciler
$ws.Range("A10").Value = 45833
$ws.Range("B10").Value = 0.875
$ws.Range("C10").Value = "Ahmet Minguzzi Grubu"
$ws.Range("D10").Value = "Ravager"
$ws.Range("E10").Value = "Kural Kesiciler"

# Row 11: 25.06.2025 22:00 - Eren Bulbul Grubu: Arakli 1961 Spor - Hubus FK
$ws.Range("A11").Value = 45833
$ws.Range("B11").Value = 0.91666666666666663
$ws.Range("C11").Value = "Eren Bülbül Grubu"
$ws.Range("D11").Value = "Araklı 1961 Spor"
$ws.Range("E11").Value = "Hubuş FK"

# Rows 12 and 13 already existed (empty, but pre-styled) - just fill values.
# Row 12: 26.06.2025 21:00 - Eren Bulbul Grubu: Armedospor - Of 1461
$ws.Range("A12").Value = 45834
$ws.Range("B12").Value = 0.875
$ws.Range("C12").Value = "Eren Bülbül Grubu"
$ws.Range("D12").Value = "Armedospor"
$ws.Range("E12").Value = "Of 1461"

# Row 13: 26.06.2025 22:00 - Narin Guran Grubu: Of FK - 61.Alay
$ws.Range("A13").Value = 45834
$ws.Range("B13").Value = 0.91666666666666663
$ws.Range("C13").Value = "Narin Güran Grubu"
$ws.Range("D13").Value = "Of FK"
$ws.Range("E13").Value = "61.Alay"

# Move the active selection to where the user ended up after entering data.
$ws.Range("E14").Select()
